# Convention change to support multi-axle vehicles:
# add two new worksheets ("Truck_Amandla_A2" and "Truck_Amandla_A3"),
# each a copy of the "Bus_Makhulu_r" sheet (same layout / formatting),
# appended at the end of the workbook, with their "Instance" label
# (cell H3) updated to match the new sheet name.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Bus_Makhulu_r")

# --- Truck_Amandla_A2 -------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheetA2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetA2.Name = "Truck_Amandla_A2"
$sheetA2.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A2"

# --- Truck_Amandla_A3 -------------------------------------------------
$sheetA2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheetA3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetA3.Name = "Truck_Amandla_A3"
$sheetA3.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A3"

# Last sheet added becomes the active / selected tab.
$sheetA3.Activate()
